$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 84, pushing the existing rows 84-92 down to 85-93
$ws.Rows("84:84").Insert()

# Populate the new row 84 with the newest weekly data point, matching
# the same (constant) columns used by the rest of the block and the
# new date/price values from the latest week.
$ws.Cells.Item(84, 1).Value = 8
$ws.Cells.Item(84, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(84, 3).Value = "Coquimbo"
$ws.Cells.Item(84, 4).Value = 44504
$ws.Cells.Item(84, 5).Value = 4
$ws.Cells.Item(84, 6).Value = 100112040
$ws.Cells.Item(84, 7).Value = "Cilantro"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 2400
$ws.Cells.Item(84, 11).Value = 1300
$ws.Cells.Item(84, 12).Value = 1500
$ws.Cells.Item(84, 13).Value = 1400
$ws.Cells.Item(84, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(84, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(84, 16).Value = 933
$ws.Cells.Item(84, 17).Value = 1.5
$ws.Cells.Item(84, 18).Value = "Hortaliza"

# Make sure the new D84 keeps the date-style formatting used by the
# rest of the column (same number format as D85:D93).
$ws.Cells.Item(84, 4).NumberFormat = $ws.Cells.Item(85, 4).NumberFormat
